# edit.ps1
# Applies the update to Base/consolidado.xlsx:
#  - Sheet "SPN": append rows 108-114 (week 9 incidents for Willian Jones' team)
#  - Sheet "ITI": mark rows 123 & 129 as "Resolvido"; append rows 135-158
#    (week 9 & 10 incidents for Emerson Simette's team)

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: writes a value into a cell as TEXT (string), even when the
# text looks like a date (dd/mm/yyyy) or other auto-converted value.
# This mirrors the source workbook, where every column other than
# "Semana" (C) and "Incidente" (F) is stored as plain text.
# -----------------------------------------------------------------
function Set-TextCell {
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# -----------------------------------------------------------------
# Helper: writes a full data row (columns A..J) for this workbook's
# layout: Setor, Responsavel, Semana, Inicio_Semana, Final_Semana,
# Incidente, Backlog, Data, Status, Coordenador
# -----------------------------------------------------------------
function Set-DataRow {
    param($Sheet, $RowData)

    $rowNum = $RowData[0]

    Set-TextCell $Sheet $rowNum 1  $RowData[1]   # A Setor
    Set-TextCell $Sheet $rowNum 2  $RowData[2]   # B Responsavel
    $Sheet.Cells.Item($rowNum, 3).Value = $RowData[3]   # C Semana (number)
    Set-TextCell $Sheet $rowNum 4  $RowData[4]   # D Inicio_Semana
    Set-TextCell $Sheet $rowNum 5  $RowData[5]   # E Final_Semana
    $Sheet.Cells.Item($rowNum, 6).Value = $RowData[6]   # F Incidente (number)
    Set-TextCell $Sheet $rowNum 7  $RowData[7]   # G Backlog
    Set-TextCell $Sheet $rowNum 8  $RowData[8]   # H Data
    Set-TextCell $Sheet $rowNum 9  $RowData[9]   # I Status
    Set-TextCell $Sheet $rowNum 10 $RowData[10]  # J Coordenador
}

# -----------------------------------------------------------------
# 1) Sheet "SPN": update status of two existing rows and append the
#    new week-9 rows (108-114)
# -----------------------------------------------------------------
$wsSPN = $wb.Worksheets.Item("SPN")

# Row data: RowNum, A, B, C, D, E, F, G, H, I, J
$sheet1NewRows = @(
    @(108, 'SPN', 'Arthur Hassuma', 9, '02/12/2024', '06/12/2024', 314457, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(109, 'SPN', 'Arthur Hassuma', 9, '02/12/2024', '06/12/2024', 315282, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(110, 'SPN', 'Arthur Hassuma', 9, '02/12/2024', '06/12/2024', 315511, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(111, 'SPN', 'Arthur Hassuma', 9, '02/12/2024', '06/12/2024', 315663, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(112, 'SPN', 'Higor Cruz', 9, '02/12/2024', '06/12/2024', 315374, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(113, 'SPN', 'Luan Pierry', 9, '02/12/2024', '06/12/2024', 315638, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones'),
    @(114, 'SPN', 'Luan Pierry', 9, '02/12/2024', '06/12/2024', 315817, '11/2024', '02/12/2024', 'Resolvido', 'Willian Jones')
)

foreach ($row in $sheet1NewRows) {
    Set-DataRow $wsSPN $row
}

# -----------------------------------------------------------------
# 2) Sheet "ITI": mark incidents 314870 (row 123) and 315005 (row 129)
#    as "Resolvido", then append the new week-9/10 rows (135-158)
# -----------------------------------------------------------------
$wsITI = $wb.Worksheets.Item("ITI")

Set-TextCell $wsITI 123 9 "Resolvido"
Set-TextCell $wsITI 129 9 "Resolvido"

$sheet2NewRows = @(
    @(135, 'ITI', 'Erick Silva', 9, '02/12/2024', '06/12/2024', 315595, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(136, 'ITI', 'Erick Silva', 9, '02/12/2024', '06/12/2024', 315683, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(137, 'ITI', 'Erick Silva', 9, '02/12/2024', '06/12/2024', 315754, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(138, 'ITI', 'Gustavo Linpiski', 9, '02/12/2024', '06/12/2024', 315377, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(139, 'ITI', 'Gustavo Linpiski', 9, '02/12/2024', '06/12/2024', 315966, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(140, 'ITI', 'Jorgenaldo Reis', 9, '02/12/2024', '06/12/2024', 315952, '11/2024', '02/12/2024', 'Pendente', 'Emerson Simette'),
    @(141, 'ITI', 'Jorgenaldo Reis', 9, '02/12/2024', '06/12/2024', 316151, '11/2024', '02/12/2024', 'Pendente', 'Emerson Simette'),
    @(142, 'ITI', 'Jorgenaldo Reis', 9, '02/12/2024', '06/12/2024', 315807, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(143, 'ITI', 'Jorgenaldo Reis', 9, '02/12/2024', '06/12/2024', 316110, '11/2024', '02/12/2024', 'Resolvido', 'Emerson Simette'),
    @(144, 'ITI', 'Jose Acevedo', 9, '02/12/2024', '06/12/2024', 315183, '11/2024', '02/12/2024', 'Pendente', 'Emerson Simette'),
    @(145, 'ITI', 'Alana Neris', 10, '09/12/2024', '13/12/2024', 316702, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(146, 'ITI', 'Alana Neris', 10, '09/12/2024', '13/12/2024', 316765, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(147, 'ITI', 'Alana Neris', 10, '09/12/2024', '13/12/2024', 316993, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(148, 'ITI', 'Edson Campos', 10, '09/12/2024', '13/12/2024', 315916, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(149, 'ITI', 'Edson Campos', 10, '09/12/2024', '13/12/2024', 315812, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(150, 'ITI', 'Erick Silva', 10, '09/12/2024', '13/12/2024', 316732, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(151, 'ITI', 'Erick Silva', 10, '09/12/2024', '13/12/2024', 316626, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(152, 'ITI', 'Erick Silva', 10, '09/12/2024', '13/12/2024', 316501, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(153, 'ITI', 'Gabriel Lopez', 10, '09/12/2024', '13/12/2024', 316114, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(154, 'ITI', 'Jacyr Popenda', 10, '09/12/2024', '13/12/2024', 316940, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(155, 'ITI', 'Jorgenaldo Reis', 10, '09/12/2024', '13/12/2024', 315817, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(156, 'ITI', 'Jose Acevedo', 10, '09/12/2024', '13/12/2024', 316763, '12/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(157, 'ITI', 'Lourival Moizés', 10, '09/12/2024', '13/12/2024', 315310, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette'),
    @(158, 'ITI', 'Lourival Moizés', 10, '09/12/2024', '13/12/2024', 315758, '11/2024', '09/12/2024', 'Pendente', 'Emerson Simette')
)

foreach ($row in $sheet2NewRows) {
    Set-DataRow $wsITI $row
}
